{"js": "// Translate the English subtitle text runs to Swahili, matching the\n// commit \"New translations ... (Swahili, Kenya)\".\nconst replacements = [\n  [\"The playful mathematicians:\", \"Wanahisabati wanaocheza:\"],\n  [\n    \"** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino\",\n    \"** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino\",\n  ],\n  [\"[Music]\", \"[Muziki]\"],\n  [\"there are two mathematicians, let's call\", \"kuna wanahisabati wawili, tupige simu\"],\n  [\"them Fil and Mike who meet each other\", \"Fil na Mike wanaokutana\"],\n  [\"again after a long time. After some\", \"tena baada ya muda mrefu. Baada ya baadhi\"],\n  [\"chatting, Phil says he has three children, then\", \"kuzungumza, Phil anasema ana watoto watatu, basi\"],\n  [\"Mike, astonished, asks: 'How old are they?' Fil,\", \"Kwa mshangao, Mike anauliza: 'Wana umri gani?' Fil,\"],\n  [\"being a playful mathematician, answers\", \"kuwa mwanahisabati mchezaji, anajibu\"],\n  [\"'You tell me! I'll give you a hint: if you\", \"'Wewe niambie! Nitakupa kidokezo: ikiwa wewe\"],\n  [\"multiply the three ages together you\", \"zidisheni enzi tatu pamoja ninyi\"],\n  [\"get 36.' Mike takes sometimes to think\", \"pata 36.' Mike huchukua wakati mwingine kufikiria\"],\n  [\"and says: 'I'm sorry Fil, but I do need\", \"na kusema: 'Samahani Fil, lakini nahitaji\"],\n  [\"another hint. So Fil tells Mike:\", \"kidokezo kingine. Kwa hivyo Fil anamwambia Mike:\"],\n  [\"'Yes, sure, here it is: if you had up to\", \"'Ndiyo, hakika, hapa ni: kama alikuwa na hadi\"],\n  [\"three ages you get the number of math\", \"miaka mitatu unapata idadi ya hesabu\"],\n  [\"papers we publish together. Do you remember it?'\", \"karatasi tunachapisha pamoja. Je, unaikumbuka?'\"],\n  [\"'Yes I do remember How many, but still\", \"'Ndio nakumbuka wangapi, lakini bado\"],\n  [\"I do not have enough information! I need\", \"Sina taarifa za kutosha! nahitaji\"],\n  [\"at least one more.' Fil says: 'Yes don't\", \"angalau moja zaidi.' Fil anasema: 'Ndiyo usifanye hivyo\"],\n  [\"worry but this is the last one:\", \"wasiwasi lakini hii ni ya mwisho:\"],\n  [\"The youngest one has blues eyes.' And\", \"Mdogo ana macho ya blues.' Na\"],\n  [\"suddenly Mike gets the answer. You\", \"ghafla Mike anapata jibu. Wewe\"],\n  [\"hear the conversation but you don't know\", \"sikia mazungumzo lakini hujui\"],\n  [\"how many papers they published together.\", \"ni karatasi ngapi walichapisha pamoja.\"],\n  [\"However, you do want to know the ages of\", \"Hata hivyo, unataka kujua umri wa\"],\n  [\"the three children. Can you figure them\", \"watoto watatu. Je, unaweza kuwahesabu\"],\n  [\"out?\", \"nje?\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the English subtitle text runs to Swahili, matching the\n# commit \"New translations ... (Swahili, Kenya)\".\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"The playful mathematicians:\", \"Wanahisabati wanaocheza:\"),\n  @(\"** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino\", \"** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino\"),\n  @(\"[Music]\", \"[Muziki]\"),\n  @(\"there are two mathematicians, let's call\", \"kuna wanahisabati wawili, tupige simu\"),\n  @(\"them Fil and Mike who meet each other\", \"Fil na Mike wanaokutana\"),\n  @(\"again after a long time. After some\", \"tena baada ya muda mrefu. Baada ya baadhi\"),\n  @(\"chatting, Phil says he has three children, then\", \"kuzungumza, Phil anasema ana watoto watatu, basi\"),\n  @(\"Mike, astonished, asks: 'How old are they?' Fil,\", \"Kwa mshangao, Mike anauliza: 'Wana umri gani?' Fil,\"),\n  @(\"being a playful mathematician, answers\", \"kuwa mwanahisabati mchezaji, anajibu\"),\n  @(\"'You tell me! I'll give you a hint: if you\", \"'Wewe niambie! Nitakupa kidokezo: ikiwa wewe\"),\n  @(\"multiply the three ages together you\", \"zidisheni enzi tatu pamoja ninyi\"),\n  @(\"get 36.' Mike takes sometimes to think\", \"pata 36.' Mike huchukua wakati mwingine kufikiria\"),\n  @(\"and says: 'I'm sorry Fil, but I do need\", \"na kusema: 'Samahani Fil, lakini nahitaji\"),\n  @(\"another hint. So Fil tells Mike:\", \"kidokezo kingine. Kwa hivyo Fil anamwambia Mike:\"),\n  @(\"'Yes, sure, here it is: if you had up to\", \"'Ndiyo, hakika, hapa ni: kama alikuwa na hadi\"),\n  @(\"three ages you get the number of math\", \"miaka mitatu unapata idadi ya hesabu\"),\n  @(\"papers we publish together. Do you remember it?'\", \"karatasi tunachapisha pamoja. Je, unaikumbuka?'\"),\n  @(\"'Yes I do remember How many, but still\", \"'Ndio nakumbuka wangapi, lakini bado\"),\n  @(\"I do not have enough information! I need\", \"Sina taarifa za kutosha! nahitaji\"),\n  @(\"at least one more.' Fil says: 'Yes don't\", \"angalau moja zaidi.' Fil anasema: 'Ndiyo usifanye hivyo\"),\n  @(\"worry but this is the last one:\", \"wasiwasi lakini hii ni ya mwisho:\"),\n  @(\"The youngest one has blues eyes.' And\", \"Mdogo ana macho ya blues.' Na\"),\n  @(\"suddenly Mike gets the answer. You\", \"ghafla Mike anapata jibu. Wewe\"),\n  @(\"hear the conversation but you don't know\", \"sikia mazungumzo lakini hujui\"),\n  @(\"how many papers they published together.\", \"ni karatasi ngapi walichapisha pamoja.\"),\n  @(\"However, you do want to know the ages of\", \"Hata hivyo, unataka kujua umri wa\"),\n  @(\"the three children. Can you figure them\", \"watoto watatu. Je, unaweza kuwahesabu\"),\n  @(\"out?\", \"nje?\")\n)\n\n# Use Find to locate each run of English text, then assign Range.Text\n# directly (instead of Find.Execute's Replace argument) so that Word's\n# smart-quote autocorrect does not mangle the straight apostrophes/quotes\n# that appear in the Swahili replacement text.\nforeach ($pair in $pairs) {\n  $searchText = $pair[0]\n  $replaceText = $pair[1]\n  $rng = $d.Content\n  while ($rng.Find.Execute($searchText)) {\n    $rng.Text = $replaceText\n  }\n}\n"}
